$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.749.63'
$ws.Range('E2').Value = '  +2.36%  '
$ws.Range('D3').Value = '1.873.31'
$ws.Range('E3').Value = '  +2.09%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.49'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3856'
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07871'
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9775'
$ws.Range('E10').Value = '  +1.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.85'
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('D12').Value = '1.812.62'
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.012'
$ws.Range('E13').Value = '  +1.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.707'
$ws.Range('E14').Value = '  +0.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06951'
$ws.Range('E15').Value = '  +1.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.46'
$ws.Range('E16').Value = '  +1.38%  '
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001003'
$ws.Range('E18').Value = '  +0.81%  '
$ws.Range('E19').Value = '  +0.63%  '
$ws.Range('D21').Value = '28.747.83'
$ws.Range('E21').Value = '  +2.29%  '
$ws.Range('E22').Value = '  -1.05%  '
$ws.Range('E23').Value = '  +0.70%  '
$ws.Range('D25').Value = '2.096.22'
$ws.Range('E25').Value = '  +4.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.00'
$ws.Range('E26').Value = '  -0.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.32'
$ws.Range('E27').Value = '  +0.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.890'
$ws.Range('E28').Value = '  +3.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.989'
$ws.Range('E29').Value = '  +1.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '119.30'
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09319'
$ws.Range('E31').Value = '  +0.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9188'
$ws.Range('E32').Value = '  -2.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.295'
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.336'
$ws.Range('E34').Value = '  +0.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.322'
$ws.Range('E35').Value = '  +0.77%  '
$ws.Range('E36').Value = '  -0.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.153'
$ws.Range('E37').Value = '  +0.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02080'
$ws.Range('E38').Value = '  -2.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.659'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5630'
$ws.Range('E40').Value = '  +0.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1783'
$ws.Range('E41').Value = '  +1.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.781'
$ws.Range('E42').Value = '  -1.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.07233'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.76'
$ws.Range('E44').Value = '  +1.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5294'
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.146'
$ws.Range('E46').Value = '  +0.44%  '
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '113.04'
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.408'
$ws.Range('E50').Value = '  +3.64%  '
$ws.Range('E51').Value = '  +0.32%  '
